$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '26.980.18'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  +2.03%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '1.841.02'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  +1.72%  '
$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '1.008'
$cell.Style = 'Normal'
$ws.Range('E4').Value = '  +0.26%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '309.35'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +1.18%  '
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('E7').Value = '  +3.75%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.3620'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  +1.04%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.07126'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +0.96%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.9107'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  +2.46%  '
$ws.Range('E11').Value = '  +0.90%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '0.07676'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -1.37%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '1.826.23'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  +0.77%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '5.266'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  +0.00%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '6.382'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  +1.24%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '88.17'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  +4.06%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '1.011'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  +0.39%  '
$ws.Range('E18').Value = '  +0.79%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '1.008'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +0.26%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '27.007.55'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +1.96%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '14.29'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +0.95%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '5.008'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +0.99%  '
$ws.Range('E23').Value = '  +1.16%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '1.932'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -0.91%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '152.43'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  +0.82%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '18.18'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +2.34%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '2.026'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -1.69%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '113.94'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  +1.62%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '4.880'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  +0.85%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '0.08857'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +2.04%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '3.203'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  +2.91%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '2.814'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  +1.79%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '0.7453'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  +0.54%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '1.168'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  +5.35%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '4.458'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  +0.44%  '
$ws.Range('E36').Value = '  +1.00%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '2.975'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +2.92%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.01935'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  +0.63%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '0.05155'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  +0.80%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.5169'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  +1.86%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '6.891'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +1.99%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '0.1510'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +0.21%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '8.101'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  +0.74%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '10.44'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +5.04%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '0.4678'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +0.12%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '1.009'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +0.51%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '100.46'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  +0.83%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '1.602'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  +1.91%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '0.06040'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  +0.99%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '64.33'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  +1.27%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '36.18'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  +0.92%  '
